$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Postal Code"
$ws.Range("B1").Value = "Floor Number"
$ws.Range("C1").Value = "Unit Number"
$ws.Range("D1").Value = "Status"
$ws.Range("E1").Value = "Flat Type"

$ws.Range("D3").Select()
